$wb = $excel.ActiveWorkbook

# Add the new "design_heat_load" worksheet as the last sheet (after the
# existing "further_parameters" sheet) and make it the active sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$ws.Name = "design_heat_load"

$ws.Range("A1").Value = "House"
$ws.Range("B1").Value = "Heat load"

$ws.Range("A2").Value = "ADS_1"
$ws.Range("B2").Value = 11000

$ws.Range("A3").Value = "ADS_10"
$ws.Range("B3").Value = 16000

$ws.Range("A4").Value = "ADS_11"
$ws.Range("B4").Value = 12000

$ws.Range("A5").Value = "ADS_12"
$ws.Range("B5").Value = 15000

$ws.Range("C3").Select() | Out-Null
